$wb = $excel.ActiveWorkbook

# --- Sheet "#export" ---
$ws3 = $wb.Worksheets.Item("#export")

$ws3.Range("A4").Value = "asdf"
$ws3.Range("B4").Value = "qwer"
$ws3.Range("C4").Value = 1
$ws3.Range("D4").Value = "01_A0_Colon_T03-2017_naive_170427_UKy_GCB_rep1-quench"
$ws3.Range("E4").Value = 289287.73343735602
$ws3.Range("F4").Value = 0
$ws3.Range("G4").Value = 8490014.3650100008
$ws3.Range("H4").Value = 0
$ws3.Range("I4").Value = 439597.55237699999
$ws3.Range("J4").Value = "NA"
$ws3.Range("K4").Value = 0
$ws3.Range("L4").Value = 0
$ws3.Range("M4").Value = 20
$ws3.Range("N4").Value = 10
$ws3.Range("O4").Value = 0.618176844244679
$ws3.Range("P4").Value = 0.255757329816374
$ws3.Range("Q4").Value = 0
$ws3.Range("R4").Value = 0
$ws3.Range("S4").Value = "Compound: name of assigned metabolite, noStd means assigment was NOT verified with standard compound"

$ws3.Range("B5").Select()

# --- Sheet "#convert" ---
$ws1 = $wb.Worksheets.Item("#convert")
$ws1.Activate()

# Change #unique=true -> #match=unique on existing row 1
$ws1.Range("D1").Value = "#match=unique"

# New block at rows 4-5 (duplicate of rows 1-2, but C4 uses #measurement.assignment.assign)
$ws1.Range("A4").Value = "#tags"
$ws1.Range("B4").Value = "#measurement.compound.value"
$ws1.Range("C4").Value = "#measurement.assignment.assign"
$ws1.Range("D4").Value = "#match=unique"
$ws1.Range("E4").Value = "#comparison=levenshtein"

$ws1.Range("B5").Value = "(S)-2-Acetolactate_Glutaric acid_Methylsuccinic acid_MP_NoStd"
$ws1.Range("C5").Value = "asdf"

# New block at rows 7-8
$ws1.Range("A7").Value = "#tags"
$ws1.Range("B7").Value = "#measurement.formula.value"
$ws1.Range("C7").Value = "#measurement.assignment.assign"
$ws1.Range("D7").Value = "#match=unique"
$ws1.Range("E7").Value = "#comparison=levenshtein"

$ws1.Range("B8").Value = "qwer"
$ws1.Range("C8").Value = "qwer"

$ws1.Columns("D:D").ColumnWidth = 12.666666666666666

$ws1.Range("C3").Select()
